# Apply the "DB Schema" edit: remove the `dessert` entry from the
# `courseType` array (Collection = recipes, Field = courseType), which was
# listed alphabetically alongside breakfast/brunch/dinner/lunch/snack.
#
# Removing the whole row lets Excel do the natural work of shifting every
# cell below it up by one row (and dropping the now-unused "dessert" shared
# string), which is exactly the cascading change shown in the diff
# (dimension B1:G48 -> B1:G47, every row after it renumbered, every shared
# string index after "dessert" shifted down by one, the sortState range for
# the cuisine list shifted up by one row, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the "dessert" row defensively (in case layout ever shifts) instead
# of hard-coding row 25, by scanning column E for the courseType array
# values underneath the "courseType" row.
$targetRow = 0
$usedRows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    $val = $ws.Cells.Item($r, 5).Value2
    if ($val -eq "dessert") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows($targetRow).Delete()
}

# Re-point the worksheet's remembered sort range for the (already
# alphabetised) cuisine list, which moved up one row along with everything
# else below the deleted row.
$sortTop = 0
$usedRows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    if ($ws.Cells.Item($r, 5).Value2 -eq "Italian") {
        $sortTop = $r
        break
    }
}
if ($sortTop -gt 0) {
    $sortBottom = $sortTop + 5
    $sortRange = $ws.Range($ws.Cells.Item($sortTop, 5), $ws.Cells.Item($sortBottom, 5))
    $srt = $ws.Sort
    $srt.SortFields.Clear()
    $srt.SortFields.Add($ws.Cells.Item($sortTop, 5)) | Out-Null
    $srt.SetRange($sortRange)
    $srt.Header = -4142
    $srt.Orientation = 1
    $srt.Apply()
}

# Restore the usual "fresh open" selection at B1 (top-left data cell).
$ws.Range("B1").Select()
